$d = $word.ActiveDocument

# 1. "... y otro para la tabla de tutores." -> "... y otros dos para la tabla de tutores."
$found1 = $d.Content.Find.Execute(
    "y otro para la tabla de tutores.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "y otros dos para la tabla de tutores.",
    2)
if (-not $found1) {
    Write-Output "WARNING: phrase 1 not found"
}

# 2. "...idea de como quiero que sea el estilo final." -> "...idea de cómo quiero que sea el estilo final."
$found2 = $d.Content.Find.Execute(
    "una idea de como quiero que sea el estilo final.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "una idea de cómo quiero que sea el estilo final.",
    2)
if (-not $found2) {
    Write-Output "WARNING: phrase 2 not found"
}

# 3. The "_GoBack" bookmark (Word's "last edit location" marker) moves from the
#    very end of the document (after "Fin de entrevista") to right after the
#    "cómo" that was just typed/corrected above.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$comoRange = $d.Content.Duplicate
$foundComo = $comoRange.Find.Execute(
    "idea de cómo",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "",
    0)
if (-not $foundComo) {
    Write-Output "WARNING: 'cómo' anchor not found for bookmark placement"
}

$bookmarkRange = $d.Range($comoRange.End, $comoRange.End)
$d.Bookmarks.Add("_GoBack", $bookmarkRange) | Out-Null
